$d = $word.ActiveDocument

# 1. Replace the existing paragraph's text (which is split across several
#    runs and wrapped by proofErr spell-check markers) with a single run
#    containing the new sentence. Find/Replace consolidates the match into
#    one run and preserves the paragraph's own formatting (the underlined
#    paragraph mark), matching the diff exactly.
$d.Content.Find.Execute("Probando git desde cero con word", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Este es mi segunda prueba", 2) | Out-Null

# 2. Insert a brand-new, plain (unformatted) paragraph above it containing
#    "Probando git desde cero con Word" split into two runs. We use
#    InsertXML with a literal OOXML snippet so the new runs do NOT inherit
#    the underline formatting that sits on the existing paragraph/run at
#    the insertion point (a plain InsertBefore/InsertParagraphBefore would
#    pick that formatting up). Because this runtime's InsertXML performs a
#    literal insert (no "smart" paragraph-merge like real Word), we supply
#    a trailing empty paragraph to force the split, then delete that helper
#    empty paragraph afterwards.
$firstPara = $d.Paragraphs(1)
$insertionPoint = $firstPara.Range.Duplicate
$insertionPoint.Collapse(1)  # wdCollapseStart

$xmlSnippet = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t xml:space="preserve">Probando git desde cero con </w:t></w:r><w:r><w:t>Word</w:t></w:r></w:p>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertionPoint.InsertXML($xmlSnippet) | Out-Null

# Remove the helper empty paragraph that was introduced to force the split
# (it sits right after the freshly inserted paragraph, before the original
# one).
$helperPara = $d.Paragraphs(2)
$helperPara.Range.Delete() | Out-Null
